$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the timestamp value in A14 (fractional day value changed)
$ws.Range("A14").Value = 45864.75029034722

# Append new row 15 with the new sensor reading
$ws.Range("A15").Value = 45864.79193502794
$ws.Range("B15").Value = 2025
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 15.33
$ws.Range("E15").Value = 82.38
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 4.5
$ws.Range("H15").Value = "E"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "19:00:23"

# Match the date/time number format & style used by the rest of column A
$ws.Range("A15").NumberFormat = $ws.Range("A14").NumberFormat
